# "Export with no is_pref and no lev distance":
# re-export the id (col B) / speaker_variant (col C) pairs for rows 2-22 in
# their new (source) order, built from an exact slug of the label instead of
# a levenshtein-matched id, and drop every is_prefered (col D) "x" marker.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '#moses-hoedende-sijne-schapen-aen-den-berghe-horeb,-spreeckt:'
$ws.Range("C2").Value = 'Moses hoedende sijne Schapen aen den Berghe Horeb, spreeckt:'
$ws.Range("D2").ClearContents()

$ws.Range("B3").Value = '#phar'
$ws.Range("C3").Value = 'Phar'
$ws.Range("D3").ClearContents()

$ws.Range("B4").Value = '#aar'
$ws.Range("C4").Value = 'Aar'
$ws.Range("D4").ClearContents()

$ws.Range("B5").Value = '#cal'
$ws.Range("C5").Value = 'Cal'
$ws.Range("D5").ClearContents()

$ws.Range("B6").Value = '#moys'
$ws.Range("C6").Value = 'Moys'
$ws.Range("D6").ClearContents()

$ws.Range("B7").Value = '#pharao-den-coningh.-rey-der-egyptenaren'
$ws.Range("C7").Value = 'Pharao den Coningh. Rey der Egyptenaren'
$ws.Range("D7").ClearContents()

$ws.Range("B8").Value = '#corach,-iosua,-ende-caleb'
$ws.Range("C8").Value = 'Corach, Iosua, ende Caleb'
$ws.Range("D8").ClearContents()

$ws.Range("B9").Value = '#pharao-den-koningh'
$ws.Range("C9").Value = 'Pharao den Koningh'
$ws.Range("D9").ClearContents()

$ws.Range("B10").Value = '#fama,-oft-''t-blasende-gherucht'
$ws.Range("C10").Value = 'Fama, oft ''t blasende gherucht'
$ws.Range("D10").ClearContents()

$ws.Range("B11").Value = '#god'
$ws.Range("C11").Value = 'God'
$ws.Range("D11").ClearContents()

$ws.Range("B12").Value = '#den-reye-der-israeliten-singhen'
$ws.Range("C12").Value = 'Den Reye der Israeliten singhen'
$ws.Range("D12").ClearContents()

$ws.Range("B13").Value = '#tiph'
$ws.Range("C13").Value = 'Tiph'
$ws.Range("D13").ClearContents()

$ws.Range("B14").Value = '#man'
$ws.Range("C14").Value = 'Man'
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = '#ios'
$ws.Range("C15").Value = 'Ios'
$ws.Range("D15").ClearContents()

$ws.Range("B16").Value = '#moyses-doet-sijn-offerhande-ende-spreect'
$ws.Range("C16").Value = 'Moyses doet sijn Offerhande ende spreect'
$ws.Range("D16").ClearContents()

$ws.Range("B17").Value = '#cor'
$ws.Range("C17").Value = 'Cor'
$ws.Range("D17").ClearContents()

$ws.Range("B18").Value = '#serax'
$ws.Range("C18").Value = 'Serax'
$ws.Range("D18").ClearContents()

$ws.Range("B19").Value = '#hymne-ofte-lof-sangh-vanden-israelijtschen-reye'
$ws.Range("C19").Value = 'Hymne ofte Lof-sangh vanden Israelijtschen Reye'
$ws.Range("D19").ClearContents()

$ws.Range("B20").Value = '#vrou'
$ws.Range("C20").Value = 'Vrou'
$ws.Range("D20").ClearContents()

$ws.Range("B21").Value = '#choor'
$ws.Range("C21").Value = 'CHOOR'
$ws.Range("D21").ClearContents()

$ws.Range("B22").Value = '#hooft'
$ws.Range("C22").Value = 'Hooft'
$ws.Range("D22").ClearContents()
